$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Poroto granado" at
# Terminal Hortofrutícola Agro Chillán. It belongs chronologically right
# before the existing row 70, so insert a fresh row there (pushing the old
# rows 70-76 down to 71-77) and populate it with the new record.
$ws.Rows("70:70").Insert()

$ws.Cells.Item(70, 1).Value = 7
$ws.Cells.Item(70, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(70, 3).Value = 'Ñuble'
$ws.Cells.Item(70, 4).Value = 44585
$ws.Cells.Item(70, 5).Value = 16
$ws.Cells.Item(70, 6).Value = 100112030
$ws.Cells.Item(70, 7).Value = 'Poroto granado'
$ws.Cells.Item(70, 8).Value = 'Sin especificar'
$ws.Cells.Item(70, 9).Value = 'Primera'
$ws.Cells.Item(70, 10).Value = 100
$ws.Cells.Item(70, 11).Value = 23000
$ws.Cells.Item(70, 12).Value = 24000
$ws.Cells.Item(70, 13).Value = 23500
$ws.Cells.Item(70, 14).Value = ([char]36) + '/saco 25 kilos'
$ws.Cells.Item(70, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(70, 16).Value = 940
$ws.Cells.Item(70, 17).Value = 25
$ws.Cells.Item(70, 18).Value = 'Hortaliza'
